$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are stored as text (numeric-looking strings).
# Force text storage (NumberFormat '@') so exact digit-for-digit strings
# (including meaningful trailing zeros) are preserved, then restore the
# default 'Normal' style so no stray formatting is left behind.
$priceUpdates = @{
    "D2" = "243.12"
    "D3" = "23.16"
    "D4" = "5.399"
    "D5" = "0.05981"
    "D6" = "3.404"
    "D7" = "6.485"
    "D8" = "0.8135"
    "D9" = "0.8980"
    "D11" = "0.07404"
    "D12" = "0.03362"
    "D13" = "0.03069"
    "D14" = "0.09332"
    "D15" = "3.861"
    "D16" = "0.001576"
    "D17" = "0.04645"
    "D18" = "0.0005940"
    "D19" = "0.006081"
    "D20" = "0.005016"
    "D21" = "0.0009809"
    "D22" = "0.00007798"
    "D23" = "0.0002900"
    "D24" = "3.614"
    "D25" = "2.163"
    "D41" = "0.006189"
    "D43" = "0.002799"
    "D44" = "0.007172"
    "D45" = "0.00005186"
    "D47" = "0.0005800"
    "D50" = "0.00002099"
    "D51" = "0.0001999"
}

foreach ($ref in $priceUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$ref]
    $cell.Style = "Normal"
}
